$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Hakim Adni): clear D4/E4
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Row 8 (Kayana): clear D8/E8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()

# Row 3 (Erlenni): set D3=1, E3="SEPEDAH"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "SEPEDAH"

# Row 6 (Syfa): set D6=1, E6="SEPEDAH"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "SEPEDAH"

# Row 9 (Samuel Jofransrael): set D9=1, E9="BARBIE KEREN"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "BARBIE KEREN"
